$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Tema" to C1
$ws.Range("C1").Value = "Tema"

# Set D12 value to 4 (hours for "Recriando a Interface do Netflix")
$ws.Range("D12").Value = 4

# Update selection/view
$ws.Range("F8").Select
